$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.907.49"
$ws.Range("E2").Value = "  +5.84%  "
$ws.Range("D3").Value = "2.232.59"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "231.69"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "61.95"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").Value = "59.33"
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "2.564.91"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "15.66"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "22.02"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "5.60"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "2.255.55"
$ws.Range("D19").Value = "41.783.56"
$ws.Range("E19").Value = "  +5.63%  "
$ws.Range("D20").Value = "72.27"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "0.0₃0898"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "249.69"
$ws.Range("E23").Value = "  +9.76%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "9.70"
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").Value = "166.74"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "19.96"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "2.64"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  +6.69%  "
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("D36").Value = "0.0636"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("D38").Value = "3.65"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("E40").Value = "  +30.08%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("D43").Value = "4.85"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "8.57"
$ws.Range("E44").Value = "  +8.54%  "
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "98.92"
$ws.Range("E47").Value = "  -3.48%  "
$ws.Range("D48").Value = "1.479.94"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "16.51"
$ws.Range("E49").Value = "  -6.52%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "52.47"
$ws.Range("E51").Value = "  +8.28%  "
